# Updated cryptos list (price + 1h volume change columns, plus a row-39/40 swap).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: Excel infers numeric types from plain-looking numeric strings (e.g. "1.000" -> 1).
# The source data stores these Price values as literal text, so a leading apostrophe forces
# a text literal, and ClearFormats() drops the resulting quote-prefix cell style again so the
# cell keeps the workbook default style (no quote-prefix / text-format artifact left behind).
function Set-TextValue($cell, $text) {
    $cell.Value = "" + $text
    $cell.ClearFormats()
}

# Row 2
$ws.Range("D2").Value = "27.697.78"
$ws.Range("E2").Value = "  -0.71%  "
# Row 3
$ws.Range("D3").Value = "1.894.98"
$ws.Range("E3").Value = "  +1.15%  "
# Row 4
Set-TextValue $ws.Range("D4") "1.000"
$ws.Range("E4").Value = "  -1.15%  "
# Row 5
Set-TextValue $ws.Range("D5") "312.88"
$ws.Range("E5").Value = "  -0.35%  "
# Row 6
$ws.Range("E6").Value = "  -1.08%  "
# Row 7
Set-TextValue $ws.Range("D7") "0.4923"
$ws.Range("E7").Value = "  +1.65%  "
# Row 8
Set-TextValue $ws.Range("D8") "0.3804"
# Row 9
Set-TextValue $ws.Range("D9") "0.07333"
$ws.Range("E9").Value = "  -0.65%  "
# Row 10
Set-TextValue $ws.Range("D10") "0.9153"
$ws.Range("E10").Value = "  -2.81%  "
# Row 11
$ws.Range("E11").Value = "  -2.45%  "
# Row 12
Set-TextValue $ws.Range("D12") "0.07676"
$ws.Range("E12").Value = "  -1.81%  "
# Row 13
$ws.Range("D13").Value = "1.919.60"
$ws.Range("E13").Value = "  +2.34%  "
# Row 14
$ws.Range("E14").Value = "  -0.37%  "
# Row 15
$ws.Range("E15").Value = "  -0.18%  "
# Row 16
$ws.Range("E16").Value = "  -0.12%  "
# Row 17
$ws.Range("E17").Value = "  -1.15%  "
# Row 18
Set-TextValue $ws.Range("D18") "0.000008786"
$ws.Range("E18").Value = "  -1.15%  "
# Row 19
$ws.Range("E19").Value = "  -1.04%  "
# Row 20
$ws.Range("D20").Value = "27.877.59"
$ws.Range("E20").Value = "  -0.12%  "
# Row 21
Set-TextValue $ws.Range("D21") "14.53"
$ws.Range("E21").Value = "  -2.26%  "
# Row 22
$ws.Range("E22").Value = "  -0.06%  "
# Row 23
$ws.Range("D23").Value = "2.146.33"
$ws.Range("E23").Value = "  +0.90%  "
# Row 24
Set-TextValue $ws.Range("D24") "10.74"
$ws.Range("E24").Value = "  -1.01%  "
# Row 25
$ws.Range("E25").Value = "  -1.82%  "
# Row 26
Set-TextValue $ws.Range("D26") "153.36"
$ws.Range("E26").Value = "  -2.36%  "
# Row 27
Set-TextValue $ws.Range("D27") "18.36"
$ws.Range("E27").Value = "  -1.21%  "
# Row 28
Set-TextValue $ws.Range("D28") "2.150"
$ws.Range("E28").Value = "  +4.48%  "
# Row 29
Set-TextValue $ws.Range("D29") "115.79"
$ws.Range("E29").Value = "  -0.22%  "
# Row 30
Set-TextValue $ws.Range("D30") "4.901"
$ws.Range("E30").Value = "  -1.93%  "
# Row 31
$ws.Range("E31").Value = "  +0.08%  "
# Row 32
Set-TextValue $ws.Range("D32") "3.191"
$ws.Range("E32").Value = "  -4.11%  "
# Row 33
Set-TextValue $ws.Range("D33") "1.219"
$ws.Range("E33").Value = "  -1.07%  "
# Row 34
Set-TextValue $ws.Range("D34") "0.7664"
$ws.Range("E34").Value = "  -0.53%  "
# Row 35
$ws.Range("E35").Value = "  -0.44%  "
# Row 36
Set-TextValue $ws.Range("D36") "0.02033"
$ws.Range("E36").Value = "  -0.92%  "
# Row 37
Set-TextValue $ws.Range("D37") "2.525"
$ws.Range("E37").Value = "  -6.96%  "
# Row 38
Set-TextValue $ws.Range("D38") "1.094"
$ws.Range("E38").Value = "  -3.50%  "
# Row 39
$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws.Range("D39") "0.05280"
$ws.Range("E39").Value = "  -1.67%  "
# Row 40
$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextValue $ws.Range("D40") "0.5474"
$ws.Range("E40").Value = "  -2.62%  "
# Row 41
Set-TextValue $ws.Range("D41") "2.980"
$ws.Range("E41").Value = "  -0.57%  "
# Row 42
Set-TextValue $ws.Range("D42") "6.915"
$ws.Range("E42").Value = "  -1.97%  "
# Row 43
Set-TextValue $ws.Range("D43") "8.540"
$ws.Range("E43").Value = "  -0.63%  "
# Row 44
$ws.Range("E44").Value = "  -1.29%  "
# Row 45
Set-TextValue $ws.Range("D45") "111.69"
$ws.Range("E45").Value = "  +5.59%  "
# Row 46
Set-TextValue $ws.Range("D46") "10.60"
$ws.Range("E46").Value = "  -0.71%  "
# Row 47
Set-TextValue $ws.Range("D47") "0.4796"
$ws.Range("E47").Value = "  -1.99%  "
# Row 48
Set-TextValue $ws.Range("D48") "1.000"
$ws.Range("E48").Value = "  -1.12%  "
# Row 49
Set-TextValue $ws.Range("D49") "1.631"
$ws.Range("E49").Value = "  -2.32%  "
# Row 50
Set-TextValue $ws.Range("D50") "67.54"
$ws.Range("E50").Value = "  -0.81%  "
# Row 51
$ws.Range("E51").Value = "  -1.10%  "
